# Round the correlation-matrix percentages in B2:P6 to 2 decimal places
# (Versao 3.3 - adicionada matriz de correlacao de criterios)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.13
$ws.Range("C2").Value = 31.1
$ws.Range("D2").Value = 5.98
$ws.Range("E2").Value = 84.48
$ws.Range("F2").Value = 50.22
$ws.Range("G2").Value = 97.59
$ws.Range("H2").Value = 91.98
$ws.Range("I2").Value = 57.41
$ws.Range("J2").Value = 41.22
$ws.Range("K2").Value = 19.67
$ws.Range("L2").Value = 5.14
$ws.Range("M2").Value = 48.74
$ws.Range("N2").Value = 15.11
$ws.Range("O2").Value = 56.52
$ws.Range("P2").Value = 96.53
$ws.Range("B3").Value = 0.08
$ws.Range("C3").Value = 35.31
$ws.Range("D3").Value = 6.35
$ws.Range("E3").Value = 83.3
$ws.Range("F3").Value = 58.17
$ws.Range("G3").Value = 95.34
$ws.Range("H3").Value = 88.97
$ws.Range("I3").Value = 52.73
$ws.Range("J3").Value = 45.16
$ws.Range("K3").Value = 18.9
$ws.Range("L3").Value = 6.3
$ws.Range("M3").Value = 57.09
$ws.Range("N3").Value = 13.51
$ws.Range("O3").Value = 56.52
$ws.Range("P3").Value = 96.9
$ws.Range("B4").Value = 0.57
$ws.Range("C4").Value = 42.66
$ws.Range("D4").Value = 8.43
$ws.Range("E4").Value = 86.91
$ws.Range("F4").Value = 54.8
$ws.Range("G4").Value = 94.07
$ws.Range("H4").Value = 91.55
$ws.Range("I4").Value = 53.78
$ws.Range("J4").Value = 42.87
$ws.Range("K4").Value = 23.76
$ws.Range("L4").Value = 10.99
$ws.Range("M4").Value = 57.63
$ws.Range("N4").Value = 17.69
$ws.Range("O4").Value = 54.03
$ws.Range("P4").Value = 97.77
$ws.Range("B5").Value = 0.26
$ws.Range("C5").Value = 41.14
$ws.Range("D5").Value = 4.84
$ws.Range("E5").Value = 80.2
$ws.Range("F5").Value = 46.85
$ws.Range("G5").Value = 97.37
$ws.Range("H5").Value = 86.37
$ws.Range("I5").Value = 59.42
$ws.Range("J5").Value = 35.59
$ws.Range("K5").Value = 14.09
$ws.Range("L5").Value = 5.1
$ws.Range("M5").Value = 54.86
$ws.Range("N5").Value = 10.74
$ws.Range("O5").Value = 59.39
$ws.Range("P5").Value = 97.34
$ws.Range("B6").Value = 0.12
$ws.Range("C6").Value = 43.09
$ws.Range("D6").Value = 12.59
$ws.Range("E6").Value = 86.86
$ws.Range("F6").Value = 53.41
$ws.Range("G6").Value = 95.47
$ws.Range("H6").Value = 88.26
$ws.Range("I6").Value = 57.98
$ws.Range("J6").Value = 42.84
$ws.Range("K6").Value = 15.79
$ws.Range("L6").Value = 3.9
$ws.Range("M6").Value = 52.62
$ws.Range("N6").Value = 12.64
$ws.Range("O6").Value = 57.27
$ws.Range("P6").Value = 98.79
